$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" strings look like plain numbers (e.g. "308.59").
# The Price/Volume columns are free-form text (other rows hold values like
# "1.641.02" that are not valid numbers), so pre-format those specific cells as
# Text first -- otherwise Excel would silently convert them to numeric values
# and lose the exact printed text (and mixed-format column).
$numFmtRange = $ws.Range("D4,D5,D6,D7,D8,D9,D10,D11,D12,D13,D14,D15,D16,D18,D19,D20,D21,D22,D25,D26,D28,D29,D30,D31,D32,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51")
foreach ($area in $numFmtRange.Areas) { $area.NumberFormat = "@" }

$ws.Range("D2").Value = "23.575.88"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "1.640.86"
$ws.Range("E3").Value = "  +4.01%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "308.59"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "0.3786"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").Value = "52.99"
$ws.Range("E8").Value = "  +6.06%  "
$ws.Range("D9").Value = "0.3686"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").Value = "1.285"
$ws.Range("E10").Value = "  +5.64%  "
$ws.Range("D11").Value = "0.08218"
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "23.26"
$ws.Range("E13").Value = "  +6.40%  "
$ws.Range("D14").Value = "6.677"
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("D15").Value = "0.00001287"
$ws.Range("E15").Value = "  +5.18%  "
$ws.Range("D16").Value = "7.484"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "1.641.48"
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").Value = "95.11"
$ws.Range("E18").Value = "  +3.19%  "
$ws.Range("D19").Value = "0.06954"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").Value = "18.45"
$ws.Range("E20").Value = "  +4.56%  "
$ws.Range("D21").Value = "6.606"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "23.550.67"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("E24").Value = "  +3.27%  "
$ws.Range("D25").Value = "3.152"
$ws.Range("E25").Value = "  +11.33%  "
$ws.Range("D26").Value = "2.408"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("D28").Value = "151.86"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "5.352"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").Value = "136.63"
$ws.Range("E30").Value = "  +4.09%  "
$ws.Range("D31").Value = "2.424"
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("D32").Value = "6.858"
$ws.Range("E32").Value = "  +4.85%  "
$ws.Range("D33").Value = "1.819.11"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").Value = "0.9778"
$ws.Range("E34").Value = "  +4.97%  "
$ws.Range("D35").Value = "0.02824"
$ws.Range("E35").Value = "  +7.35%  "
$ws.Range("D36").Value = "10.50"
$ws.Range("E36").Value = "  +5.65%  "
$ws.Range("D37").Value = "0.07483"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").Value = "6.254"
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("D39").Value = "0.2550"
$ws.Range("E39").Value = "  +3.68%  "
$ws.Range("D40").Value = "0.08873"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").Value = "1.401"
$ws.Range("E41").Value = "  +4.31%  "
$ws.Range("D42").Value = "0.7196"
$ws.Range("E42").Value = "  +4.87%  "
$ws.Range("D43").Value = "12.67"
$ws.Range("E43").Value = "  +6.32%  "
$ws.Range("D44").Value = "16.17"
$ws.Range("E44").Value = "  +10.22%  "
$ws.Range("D45").Value = "0.6654"
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").Value = "2.372"
$ws.Range("E46").Value = "  +5.79%  "
$ws.Range("D47").Value = "4.048"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").Value = "0.9984"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "0.08076"
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").Value = "131.66"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").Value = "1.221"
$ws.Range("E51").Value = "  +3.41%  "

# Put those cells back to the default (General) cell style now that the text
# is safely stored, so no extra number-format is left behind on them.
foreach ($area in $numFmtRange.Areas) { $area.Style = "Normal" }
